$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 29324.914
$ws.Range("J17").Value = 27217.814
$ws.Range("L17").Value = 81653.442
$ws.Range("N17").Value = -81989.442
$ws.Range("H18").Value = 5563826
$ws.Range("I18").Value = 9259709
$ws.Range("K18").Value = 9259709
$ws.Range("M18").Value = -9259425
$ws.Range("H55").Value = 1185.4166
$ws.Range("I55").Value = 1866.6666
$ws.Range("J55").Value = 504.16666
$ws.Range("K55").Value = 1866.6666
$ws.Range("L55").Value = 504.16666
$ws.Range("M55").Value = -1652.6666
$ws.Range("N55").Value = -932.16666
$ws.Range("H70").Value = 2576
$ws.Range("I70").Value = 984.2857
$ws.Range("J70").Value = 3433.077
$ws.Range("K70").Value = 2952.8571
$ws.Range("L70").Value = 10299.231
$ws.Range("M70").Value = -2682.8571
$ws.Range("N70").Value = -10839.231
$ws.Range("H73").Value = 2576
$ws.Range("I73").Value = 984.2857
$ws.Range("J73").Value = 3433.077
$ws.Range("K73").Value = 2952.8571
$ws.Range("L73").Value = 10299.231
$ws.Range("M73").Value = -2016.8571
$ws.Range("N73").Value = -12171.231
$ws.Range("H74").Value = 3925.4736
$ws.Range("I74").Value = 3715
$ws.Range("J74").Value = 4022.6155
$ws.Range("K74").Value = 3715
$ws.Range("L74").Value = 4022.6155
$ws.Range("M74").Value = -2779
$ws.Range("N74").Value = -5894.6155
$ws.Range("H77").Value = 3925.4736
$ws.Range("I77").Value = 3715
$ws.Range("J77").Value = 4022.6155
$ws.Range("K77").Value = 18575
$ws.Range("L77").Value = 20113.0775
$ws.Range("M77").Value = -13895
$ws.Range("N77").Value = -29473.0775
$ws.Range("H121").Value = 1266.6666
$ws.Range("J121").Value = 3000
$ws.Range("L121").Value = 9000
$ws.Range("N121").Value = -12494
$ws.Range("H129").Value = 3126290.8
$ws.Range("I129").Value = 13889500
$ws.Range("J129").Value = 1488.1936
$ws.Range("K129").Value = 41668500
$ws.Range("L129").Value = 4464.5808
$ws.Range("M129").Value = -41663500
$ws.Range("N129").Value = -14464.5808
$ws.Range("H137").Value = 4171366.5
$ws.Range("I137").Value = 6255206
$ws.Range("J137").Value = 3687.5
$ws.Range("K137").Value = 18765618
$ws.Range("L137").Value = 11062.5
$ws.Range("M137").Value = -18763068
$ws.Range("N137").Value = -16162.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2381.6304
$ws.Range("I61").Value = 1357.24
$ws.Range("J61").Value = 3601.1428
$ws.Range("K61").Value = 1357.24
$ws.Range("L61").Value = 3601.1428
$ws.Range("M61").Value = -1145.24
$ws.Range("N61").Value = -4025.1428
$ws.Range("H74").Value = 734.0769
$ws.Range("I74").Value = 654.3
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 654.3
$ws.Range("L74").Value = 1000
$ws.Range("M74").Value = 219.7
$ws.Range("N74").Value = -2748
$ws.Range("H77").Value = 734.0769
$ws.Range("I77").Value = 654.3
$ws.Range("J77").Value = 1000
$ws.Range("K77").Value = 3271.5
$ws.Range("L77").Value = 5000
$ws.Range("M77").Value = 1096.5
$ws.Range("N77").Value = -13736
$ws.Range("H132").Value = 35717588
$ws.Range("I132").Value = 58826156
$ws.Range("J132").Value = 4346.5454
$ws.Range("K132").Value = 176478468
$ws.Range("L132").Value = 13039.6362
$ws.Range("M132").Value = -176475938
$ws.Range("N132").Value = -18099.6362
$ws.Range("H136").Value = 2381.6304
$ws.Range("I136").Value = 1357.24
$ws.Range("J136").Value = 3601.1428
$ws.Range("K136").Value = 4071.72
$ws.Range("L136").Value = 10803.4284
$ws.Range("M136").Value = -1521.72
$ws.Range("N136").Value = -15903.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2146.56
$ws.Range("I134").Value = 1687.0588
$ws.Range("J134").Value = 3123
$ws.Range("K134").Value = 5061.1764
$ws.Range("L134").Value = 9369
$ws.Range("M134").Value = -2526.1764
$ws.Range("N134").Value = -14439

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2276106.5
$ws.Range("I31").Value = 2633939.2
$ws.Range("K31").Value = 2633939.2
$ws.Range("M31").Value = -2633644.2
$ws.Range("H34").Value = 2276106.5
$ws.Range("I34").Value = 2633939.2
$ws.Range("K34").Value = 2633939.2
$ws.Range("M34").Value = -2633737.2
$ws.Range("H58").Value = 12502379
$ws.Range("I58").Value = 1310.3914
$ws.Range("J58").Value = 29415590
$ws.Range("K58").Value = 1310.3914
$ws.Range("L58").Value = 29415590
$ws.Range("M58").Value = -1107.3914
$ws.Range("N58").Value = -29415996
$ws.Range("H94").Value = 26317854
$ws.Range("I94").Value = 2118.3333
$ws.Range("J94").Value = 38463576
$ws.Range("K94").Value = 2118.3333
$ws.Range("L94").Value = 38463576
$ws.Range("M94").Value = -1667.3333
$ws.Range("N94").Value = -38464478
$ws.Range("H132").Value = 2829.9
$ws.Range("I132").Value = 1973.7826
$ws.Range("K132").Value = 5921.3478
$ws.Range("M132").Value = -3391.3478
$ws.Range("H134").Value = 1410.6923
$ws.Range("I134").Value = 876.3043
$ws.Range("J134").Value = 2178.875
$ws.Range("K134").Value = 2628.9129
$ws.Range("L134").Value = 6536.625
$ws.Range("M134").Value = -93.91290000000026
$ws.Range("N134").Value = -11606.625
$ws.Range("H136").Value = 12502379
$ws.Range("I136").Value = 1310.3914
$ws.Range("J136").Value = 29415590
$ws.Range("K136").Value = 3931.1742
$ws.Range("L136").Value = 88246770
$ws.Range("M136").Value = -1381.1742
$ws.Range("N136").Value = -88251870

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 4266.6665
$ws.Range("J43").Value = 4266.6665
$ws.Range("L43").Value = 12799.9995
$ws.Range("N43").Value = -13027.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H87").Value = 30000
$ws.Range("J87").Value = 30000
$ws.Range("L87").Value = 30000
$ws.Range("N87").Value = -32496
$ws.Range("H90").Value = 30000
$ws.Range("J90").Value = 30000
$ws.Range("L90").Value = 90000
$ws.Range("N90").Value = -102480
$ws.Range("H132").Value = 3124.1428
$ws.Range("I132").Value = 2732.84
$ws.Range("J132").Value = 3699.5881
$ws.Range("K132").Value = 8198.52
$ws.Range("L132").Value = 11098.7643
$ws.Range("M132").Value = -5668.52
$ws.Range("N132").Value = -16158.7643

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3024.5
$ws.Range("I7").Value = 1999
$ws.Range("J7").Value = 4050
$ws.Range("K7").Value = 1999
$ws.Range("L7").Value = 4050
$ws.Range("M7").Value = -1887
$ws.Range("N7").Value = -4274
$ws.Range("H68").Value = 1718.2916
$ws.Range("I68").Value = 1010.86365
$ws.Range("J68").Value = 9500
$ws.Range("K68").Value = 1010.86365
$ws.Range("L68").Value = 9500
$ws.Range("M68").Value = -261.86365
$ws.Range("N68").Value = -10998
$ws.Range("H71").Value = 1718.2916
$ws.Range("I71").Value = 1010.86365
$ws.Range("J71").Value = 9500
$ws.Range("K71").Value = 5054.31825
$ws.Range("L71").Value = 47500
$ws.Range("M71").Value = -1310.31825
$ws.Range("N71").Value = -54988
$ws.Range("H126").Value = 3024.5
$ws.Range("I126").Value = 1999
$ws.Range("J126").Value = 4050
$ws.Range("K126").Value = 5997
$ws.Range("L126").Value = 12150
$ws.Range("M126").Value = -3527
$ws.Range("N126").Value = -17090
$ws.Range("H132").Value = 3352.3572
$ws.Range("I132").Value = 2176.2856
$ws.Range("J132").Value = 4528.4287
$ws.Range("K132").Value = 6528.8568
$ws.Range("L132").Value = 13585.2861
$ws.Range("M132").Value = -3998.8568
$ws.Range("N132").Value = -18645.2861
$ws.Range("H136").Value = 3228408.5
$ws.Range("I136").Value = 5884727
$ws.Range("J136").Value = 2878.5
$ws.Range("K136").Value = 17654181
$ws.Range("L136").Value = 8635.5
$ws.Range("M136").Value = -17651631
$ws.Range("N136").Value = -13735.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1671.8572
$ws.Range("I62").Value = 1460.8
$ws.Range("J62").Value = 2199.5
$ws.Range("K62").Value = 1460.8
$ws.Range("L62").Value = 2199.5
$ws.Range("M62").Value = -836.8
$ws.Range("N62").Value = -3447.5
$ws.Range("H65").Value = 1671.8572
$ws.Range("I65").Value = 1460.8
$ws.Range("J65").Value = 2199.5
$ws.Range("K65").Value = 7304
$ws.Range("L65").Value = 10997.5
$ws.Range("M65").Value = -4184
$ws.Range("N65").Value = -17237.5
$ws.Range("H69").Value = 29875
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 29875
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 29875
$ws.Range("M69").Value = ""
$ws.Range("N69").Value = -31373
$ws.Range("H72").Value = 29875
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 29875
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 89625
$ws.Range("M72").Value = ""
$ws.Range("N72").Value = -97113
$ws.Range("H116").Value = 32500
$ws.Range("J116").Value = 32500
$ws.Range("L116").Value = 32500
$ws.Range("N116").Value = -41678
$ws.Range("H132").Value = 218631.62
$ws.Range("I132").Value = 279258.78
$ws.Range("J132").Value = 36750.168
$ws.Range("K132").Value = 837776.3400000001
$ws.Range("L132").Value = 110250.504
$ws.Range("M132").Value = -835246.3400000001
$ws.Range("N132").Value = -115310.504
$ws.Range("H136").Value = 1024.6586
$ws.Range("I136").Value = 644
$ws.Range("J136").Value = 1944.5834
$ws.Range("K136").Value = 1932
$ws.Range("L136").Value = 5833.7502
$ws.Range("M136").Value = 618
$ws.Range("N136").Value = -10933.7502
